$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 15, shifting existing rows 15-124 down to 16-125.
$ws.Rows(15).Insert()

# Populate the newly inserted row 15 with the new weekly data entry.
$ws.Range("A15").Value = 9
$ws.Range("B15").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C15").Value = "Metropolitana"
$ws.Range("D15").Value = "2023-02-02"
$ws.Range("E15").Value = 13
$ws.Range("F15").Value = "Fruta"
$ws.Range("G15").Value = 100101
$ws.Range("H15").Value = "Berries"
$ws.Range("I15").Value = 100101004
$ws.Range("J15").Value = "Frambuesa"
$ws.Range("K15").Value = "Sin especificar"
$ws.Range("L15").Value = "Primera"
$ws.Range("M15").Value = 350
$ws.Range("N15").Value = 8000
$ws.Range("O15").Value = 8000
$ws.Range("P15").Value = 8000
$ws.Range("Q15").Value = "$/bandeja 2 kilos"
$ws.Range("R15").Value = "Provincia de Curicó"
$ws.Range("S15").Value = 4000
$ws.Range("T15").Value = 2
